$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a handful of existing values in rows 9, 11, 12, 13, 22 ---
$ws.Range("N9").Value = 9399.04
$ws.Range("O9").Value = 9157.96

$ws.Range("N11").Value = 9511.36

$ws.Range("N12").Value = 417328.94
$ws.Range("O12").Value = 355865.13

$ws.Range("N13").Value = 93156.52
$ws.Range("O13").Value = 93156.52

$ws.Range("K22").Value = 162957.44

# --- Unmerge all the merged ranges that live inside the block that is
# about to be reshaped (rows 8-28), so the row insertion below does not
# apply its own (inconsistent) automatic merge-resizing logic. We will
# recreate the correct merges explicitly afterwards. ---
$ws.Range("A8:A23").UnMerge()
$ws.Range("B8:B23").UnMerge()
$ws.Range("C8:C18").UnMerge()
$ws.Range("D8:D18").UnMerge()
$ws.Range("E8:E16").UnMerge()
$ws.Range("F8:F16").UnMerge()
$ws.Range("G8:G16").UnMerge()
$ws.Range("H8:H16").UnMerge()
$ws.Range("E17:E18").UnMerge()
$ws.Range("F17:F18").UnMerge()
$ws.Range("G17:G18").UnMerge()
$ws.Range("H17:H18").UnMerge()
$ws.Range("C19:C23").UnMerge()
$ws.Range("D19:D23").UnMerge()
$ws.Range("E19:E21").UnMerge()
$ws.Range("F19:F21").UnMerge()
$ws.Range("G19:G21").UnMerge()
$ws.Range("H19:H21").UnMerge()
$ws.Range("E22:E23").UnMerge()
$ws.Range("F22:F23").UnMerge()
$ws.Range("G22:G23").UnMerge()
$ws.Range("H22:H23").UnMerge()
$ws.Range("A24:A28").UnMerge()
$ws.Range("B24:B28").UnMerge()
$ws.Range("C24:C25").UnMerge()
$ws.Range("D24:D25").UnMerge()
$ws.Range("E24:E25").UnMerge()
$ws.Range("F24:F25").UnMerge()
$ws.Range("G24:G25").UnMerge()
$ws.Range("H24:H25").UnMerge()
$ws.Range("C26:C28").UnMerge()
$ws.Range("D26:D28").UnMerge()
$ws.Range("E26:E27").UnMerge()
$ws.Range("F26:F27").UnMerge()
$ws.Range("G26:G27").UnMerge()
$ws.Range("H26:H27").UnMerge()

# --- Insert a new row before the current row 24 ---
# (this shifts the old rows 24-28 down to 25-29)
$ws.Rows("24:24").Insert()

# Copy the formatting of the (now unchanged) row 23 into the new row 24
# so the new row uses the same style set (s="10/11/12/13") as the rest
# of the data block, rather than Excel's auto-generated blended style.
$ws.Range("A23:O23").Copy()
$ws.Range("A24:O24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row 24 with its specific content.
# I24 ("339037") looks like a number, so prefix with an apostrophe to force
# it to be stored as text, matching the rest of the numeric-looking codes
# in this sheet (e.g. "152420", "339040", ...) which are text, not numbers.
$ws.Range("I24").Value = "'339037"
$ws.Range("J24").Value = "LOCACAO DE MAO-DE-OBRA"
$ws.Range("K24").Value = 103873.55

# The apostrophe-prefix trick above marks the cell with a "quote prefix"
# which Excel implements as a different cell style (s="11" + quotePrefix).
# Re-paste the formatting from the identically-styled I23 cell (which is
# plain style 11, no quote prefix) to restore the original style index.
$ws.Range("I23").Copy()
$ws.Range("I24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the value that changed within the shifted block (old row 27 -> new row 28) ---
$ws.Range("N28").Value = 63495

# --- Recreate the merges with their final (post-insert) extents ---
$ws.Range("A8:A24").Merge()
$ws.Range("B8:B24").Merge()
$ws.Range("C8:C18").Merge()
$ws.Range("D8:D18").Merge()
$ws.Range("E8:E16").Merge()
$ws.Range("F8:F16").Merge()
$ws.Range("G8:G16").Merge()
$ws.Range("H8:H16").Merge()
$ws.Range("E17:E18").Merge()
$ws.Range("F17:F18").Merge()
$ws.Range("G17:G18").Merge()
$ws.Range("H17:H18").Merge()
$ws.Range("C19:C24").Merge()
$ws.Range("D19:D24").Merge()
$ws.Range("E19:E21").Merge()
$ws.Range("F19:F21").Merge()
$ws.Range("G19:G21").Merge()
$ws.Range("H19:H21").Merge()
$ws.Range("E22:E24").Merge()
$ws.Range("F22:F24").Merge()
$ws.Range("G22:G24").Merge()
$ws.Range("H22:H24").Merge()
$ws.Range("A25:A29").Merge()
$ws.Range("B25:B29").Merge()
$ws.Range("C25:C26").Merge()
$ws.Range("D25:D26").Merge()
$ws.Range("E25:E26").Merge()
$ws.Range("F25:F26").Merge()
$ws.Range("G25:G26").Merge()
$ws.Range("H25:H26").Merge()
$ws.Range("C27:C29").Merge()
$ws.Range("D27:D29").Merge()
$ws.Range("E27:E28").Merge()
$ws.Range("F27:F28").Merge()
$ws.Range("G27:G28").Merge()
$ws.Range("H27:H28").Merge()
